$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = -0.5718192715580426
$ws.Range("B5").Value = 0.02513501597327435
$ws.Range("B6").Value = -0.6777052466944516
$ws.Range("B7").Value = -0.04687708722170142
$ws.Range("B8").Value = -0.6491824449565908
$ws.Range("B9").Value = -0.03143360567770964
$ws.Range("B10").Value = 0.3822062224503651
$ws.Range("B11").Value = -0.04722993886919886
$ws.Range("B12").Value = -0.5591348797726128
$ws.Range("B13").Value = -0.3824674568541571
$ws.Range("B14").Value = 0.3
$ws.Range("B15").Value = 0
$ws.Range("B16").Value = 0
$ws.Range("B17").Value = -0.07670146240877694
$ws.Range("B18").Value = -0.2
$ws.Range("B19").Value = 0.1
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = -0.1
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0.3
$ws.Range("B25").Value = -0.3
$ws.Range("B26").Value = -0.2
$ws.Range("B28").Value = 0.07231168587756959
$ws.Range("B29").Value = -0.06367835427788604
$ws.Range("B30").Value = 0.08504296140338079
$ws.Range("B31").Value = 0.0003037193523987569
$ws.Range("B32").Value = -0.01801114345337739
$ws.Range("B33").Value = 0.03826442504165321
$ws.Range("B34").Value = -0.02232391446496779
$ws.Range("B35").Value = -0.06208406292622634
$ws.Range("B36").Value = 0.006968764489216402
$ws.Range("B37").Value = -0.121948034731228
$ws.Range("B38").Value = -0.0361595126957486
$ws.Range("B39").Value = -0.08573317974146363
$ws.Range("B40").Value = 0.004670068289308601
$ws.Range("B41").Value = 0.04099992594395328
$ws.Range("B42").Value = 0.05927521680271484
$ws.Range("B43").Value = 0.04061467511413108
$ws.Range("B44").Value = 0.08363338226171732
$ws.Range("B45").Value = -0.184542838621156
$ws.Range("B46").Value = -0.4
$ws.Range("B47").Value = -0.3200594507515429
$ws.Range("B48").Value = 0.1
$ws.Range("B49").Value = -0.2628511487790233
$ws.Range("B50").Value = -0.3133720122661878
$ws.Range("B51").Value = 0.5469624758293199
$ws.Range("B52").Value = 1.1
$ws.Range("B53").Value = -0.4440571223929872
$ws.Range("B54").Value = -0.7255945204468831
$ws.Range("B55").Value = -0.5292660609007143
$ws.Range("B56").Value = -0.1550786956675604
$ws.Range("B57").Value = -2.168330733759602
$ws.Range("B58").Value = -0.03982694963614287
$ws.Range("B59").Value = 0.2669401745841223
$ws.Range("B60").Value = 0.03791487406588956
$ws.Range("B61").Value = -0.04567208272808071
$ws.Range("B62").Value = -0.5154625125417773
$ws.Range("B63").Value = -0.1813602613933202
$ws.Range("B64").Value = -0.01480819732384536
$ws.Range("B65").Value = 0.02918400950819283
$ws.Range("B66").Value = -0.03321544329283629
$ws.Range("B67").Value = 0.00001303303454188581
$ws.Range("B68").Value = -0.006125572440376981
$ws.Range("B69").Value = 0.04879937325030748
$ws.Range("B70").Value = 0.0477695913607396
$ws.Range("B71").Value = 0.4714513528429705
$ws.Range("B72").Value = -0.02605454389395597
$ws.Range("B73").Value = 0.04549112474043772

$ws.Range("A74:B82").EntireRow.Delete()
